# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.840.29'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.19%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.113.26'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.84%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.39%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.43'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.38%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.109.53'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.79%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.522'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.54%  '

$ws.Range("E10").Value = '  -2.92%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.153'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.97%  '

$ws.Range("E12").Value = '  +0.07%  '

$ws.Range("E13").Value = '  -1.45%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.20'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.48%  '

$ws.Range("E15").Value = '  -1.05%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.627.92'
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.780.09'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.09%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.16'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.92%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.114.08'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.73%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.36'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.77%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '476.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.65%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.714'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.07%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.95'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.28%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.39'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.38%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.07'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.06%  '

$ws.Range("E26").Value = '  -2.46%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.12'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.20%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.01%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.88'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.89%  '

$ws.Range("E30").Value = '  -2.39%  '

$ws.Range("E31").Value = '  -0.02%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.57'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.35%  '

$ws.Range("E33").Value = '  +0.30%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0941'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -8.40%  '

$ws.Range("E35").Value = '  -0.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.86'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.53%  '

$ws.Range("E37").Value = '  -2.92%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '46.88'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.73%  '

$ws.Range("E39").Value = '  -3.26%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.06'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.39%  '

$ws.Range("E41").Value = '  -1.77%  '

$ws.Range("E42").Value = '  -0.49%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.73'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.54%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.830.07'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.31%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '384.86'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.58%  '

$ws.Range("E46").Value = '  -1.50%  '

$ws.Range("E47").Value = '  -8.78%  '

$ws.Range("E48").Value = '  +0.82%  '

$ws.Range("E49").Value = '  -0.01%  '

$ws.Range("E51").Value = '  -1.89%  '

